# Section 3: Entering and Editing Text and Formulas Quiz
# Adds a "Percent" column (F) to the monthly budget sheet:
#   - F3 header "Percent"
#   - E4:E8 row totals (Bills+Rent, etc. across the three months)
#   - F4:F8 each row's share of the grand total
#   - B9:D9 column totals, E9 grand total, F9 = 100%

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F3").Value = "Percent"

# Per-category row totals and percentage-of-total
$rows = 4, 5, 6, 7, 8
foreach ($r in $rows) {
    $ws.Range("E$r").Formula = "=B$r+C$r+D$r"
    $ws.Range("F$r").Formula = '=E' + $r + '/$E$9'
}

# Totals row
$ws.Range("B9").Formula = "=B4+B5+B6+B7+B8"
$ws.Range("C9").Formula = "=C4+C5+C6+C7+C8"
$ws.Range("D9").Formula = "=D4+D5+D6+D7+D8"
$ws.Range("E9").Formula = "=B9+C9+D9"
$ws.Range("F9").Formula = '=E9/$E$9'

# Matches the saved selection state in the target workbook
$ws.Range("H14").Select() | Out-Null
